$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.523.87"
$ws.Range("E2").Value = "  +4.04%  "
$ws.Range("D3").Value = "2.271.36"
$ws.Range("E3").Value = "  +2.84%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'321.06"
$ws.Range("E5").Value = "  +1.78%  "
$ws.Range("D6").Value = "'105.39"
$ws.Range("E6").Value = "  +6.34%  "
$ws.Range("D7").Value = "'0.592"
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.572"
$ws.Range("E9").Value = "  +2.19%  "
$ws.Range("D10").Value = "'38.56"
$ws.Range("E10").Value = "  +4.64%  "
$ws.Range("D11").Value = "'0.0843"
$ws.Range("E11").Value = "  +1.95%  "
$ws.Range("D12").Value = "'7.87"
$ws.Range("E12").Value = "  +2.13%  "
$ws.Range("E13").Value = "  +0.71%  "
$ws.Range("D14").Value = "'0.883"
$ws.Range("E14").Value = "  +3.02%  "
$ws.Range("D15").Value = "2.620.60"
$ws.Range("E15").Value = "  +2.88%  "
$ws.Range("D16").Value = "'14.58"
$ws.Range("E16").Value = "  +3.06%  "
$ws.Range("D17").Value = "2.272.98"
$ws.Range("E17").Value = "  +2.87%  "
$ws.Range("D18").Value = "44.404.56"
$ws.Range("E18").Value = "  +3.98%  "
$ws.Range("D19").Value = "'13.90"
$ws.Range("E19").Value = "  -4.43%  "
$ws.Range("D20").Value = "'0.0000100"
$ws.Range("E20").Value = "  +4.61%  "
$ws.Range("D21").Value = "'6.53"
$ws.Range("E21").Value = "  +2.17%  "
$ws.Range("D22").Value = "'66.41"
$ws.Range("E22").Value = "  +1.82%  "
$ws.Range("E23").Value = "  +2.17%  "
$ws.Range("D24").Value = "'239.51"
$ws.Range("E24").Value = "  +1.63%  "
$ws.Range("D25").Value = "'2.21"
$ws.Range("E25").Value = "  +4.07%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'10.17"
$ws.Range("E27").Value = "  +1.79%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.22"
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "'38.39"
$ws.Range("E29").Value = "  +12.35%  "
$ws.Range("D30").Value = "'6.48"
$ws.Range("E30").Value = "  +2.88%  "
$ws.Range("D31").Value = "'20.64"
$ws.Range("E31").Value = "  +0.71%  "
$ws.Range("D32").Value = "'0.0885"
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("D33").Value = "'161.11"
$ws.Range("E33").Value = "  +3.54%  "
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("E35").Value = "  +10.40%  "
$ws.Range("E36").Value = "  +4.43%  "
$ws.Range("D37").Value = "'3.14"
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("D38").Value = "'0.121"
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("D39").Value = "'3.93"
$ws.Range("E39").Value = "  +1.80%  "
$ws.Range("D40").Value = "'4.45"
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0329"
$ws.Range("E41").Value = "  +1.51%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "'15.55"
$ws.Range("E42").Value = "  +24.61%  "
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").Value = "1.770.48"
$ws.Range("E44").Value = "  -6.31%  "
$ws.Range("D45").Value = "'0.209"
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("D46").Value = "'86.68"
$ws.Range("E46").Value = "  -1.52%  "
$ws.Range("E47").Value = "  +2.14%  "
$ws.Range("D48").Value = "'60.50"
$ws.Range("E48").Value = "  -0.96%  "
$ws.Range("D49").Value = "'75.45"
$ws.Range("E49").Value = "  -1.18%  "
$ws.Range("E50").Value = "  +7.32%  "
$ws.Range("D51").Value = "'104.27"
$ws.Range("E51").Value = "  +1.98%  "
